# Site updated: 2022-10-11 16:27:13
# Adds four new game rows (Good Job, The Lightbringer, Picross S1-S8, Zelda)
# into the already pinyin-sorted game list, re-pointing hyperlinks and
# refreshing the sheet's sort/selection metadata to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the new rows by inserting blank rows at the positions
#    that the sorted list needs (pinyin order: Good Job(gan), Lightbringer
#    (guang), Picross(hui) come before the existing Big Brain(ling) row,
#    and Zelda(sai) lands between Big Brain and Digimon(shu)).
# ---------------------------------------------------------------------
$ws.Rows(5).Insert()
$ws.Rows(5).Insert()
$ws.Rows(5).Insert()
$ws.Rows(9).Insert()

# The Insert() calls above copy the row-4 formatting (styles 5/6) onto the
# fresh rows; reset the touched columns back to the un-styled / hyperlink
# styles the new rows actually need before filling in values. (Row 7 has
# no D-column value, so its D cell is left completely untouched.)
$ws.Range("A5:A7").Style = "常规"
$ws.Range("D5:D6").Style = "常规"
$ws.Range("G5:G7").Style = "超链接"
$ws.Range("A9").Style = "常规"
$ws.Range("D9").Style = "常规"
$ws.Range("G9").Style = "超链接"

# ---------------------------------------------------------------------
# 2. Fill in the new row contents.
# ---------------------------------------------------------------------
# Row 5: Good Job
$ws.Cells.Item(5, 1).Value = "干得漂亮 | Good Job"
$ws.Cells.Item(5, 4).Value = "h4qm"
$ws.Cells.Item(5, 6).Value = "switch《干得漂亮 Good Job》xci汉化整合版下载"
$ws.Cells.Item(5, 7).Value = "https://pan.baidu.com/s/1cIIZZv89eBKv255fKtD4cQ"

# Row 6: The Lightbringer
$ws.Cells.Item(6, 1).Value = "光明使者 | The Lightbringer"
$ws.Cells.Item(6, 2).Value = 1.2
$ws.Cells.Item(6, 4).Value = "hzoy"
$ws.Cells.Item(6, 5).Value = "游戏年轮bibgame.com发布"
$ws.Cells.Item(6, 6).Value = "switch《光明使者 The Lightbringer》中文版nsp/xci下载【含1.2补丁】"
$ws.Cells.Item(6, 7).Value = "https://pan.baidu.com/s/1zGN-60z0e43PEJYohxM8cA"

# Row 7: Picross S1-S8
$ws.Cells.Item(7, 1).Value = "绘图方块系列合集 | Picross S1-S8"
$ws.Cells.Item(7, 6).Value = "switch《绘图方块系列合集》Picross S1-S8+世嘉版nsp下载【含最新补丁】"
$ws.Cells.Item(7, 7).Value = "https://pan.baidu.com/s/1EWpU8lPT_bYMm3uoK174kA?pwd=5vch"

# Row 9: Zelda: Breath of the Wild
$ws.Cells.Item(9, 1).Value = "塞尔达传说荒野之息"
$ws.Cells.Item(9, 2).Value = 1.6
$ws.Cells.Item(9, 3).Value = "10.1.1"
$ws.Cells.Item(9, 4).Value = "ph5q"
$ws.Cells.Item(9, 5).Value = "bibgame.com"
$ws.Cells.Item(9, 6).Value = "switch《塞尔达传说荒野之息》本体+v1.6整合"
$ws.Cells.Item(9, 7).Value = "https://pan.baidu.com/s/1sF9Suvp0mVQNP56NZmL64A"

# ---------------------------------------------------------------------
# 3. Hyperlinks don't follow the row-insert shifts automatically, so
#    rebuild the full hyperlink set against the final row numbers.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G3"), "https://pan.baidu.com/s/1x_V0cQZyzhAIzr97GCznlA?pwd=5fub", "list/path=%2F", "", "https://pan.baidu.com/s/1x_V0cQZyzhAIzr97GCznlA?pwd=5fub - list/path=%2F")
$ws.Hyperlinks.Add($ws.Range("G11"), "https://pan.baidu.com/s/1ZbWggC3GDJv7BUgxTIbGzg")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://pan.baidu.com/s/1KykYnfqctZOEDgJp_nxGsA?pwd=uqer")
$ws.Hyperlinks.Add($ws.Range("G10"), "https://pan.baidu.com/s/1zKgW1pjqUnZ2dtEq2xXFMw?pwd=ccx6")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://pan.baidu.com/s/1UARljz8BQP1uTU3Lie_2oQ")
$ws.Hyperlinks.Add($ws.Range("G13"), "https://pan.baidu.com/s/1Re4OiBosRO_y77sDJRBRuw", "list/path=%2F", "", "https://pan.baidu.com/s/1Re4OiBosRO_y77sDJRBRuw - list/path=%2F")
$ws.Hyperlinks.Add($ws.Range("G8"), "https://pan.baidu.com/share/init?surl=n6ivaYdevwiyNpXc1Fgpxg")
$ws.Hyperlinks.Add($ws.Range("G6"), "https://pan.baidu.com/s/1zGN-60z0e43PEJYohxM8cA")
$ws.Hyperlinks.Add($ws.Range("G9"), "https://pan.baidu.com/s/1sF9Suvp0mVQNP56NZmL64A")
$ws.Hyperlinks.Add($ws.Range("G7"), "https://pan.baidu.com/s/1EWpU8lPT_bYMm3uoK174kA?pwd=5vch")
$ws.Hyperlinks.Add($ws.Range("G5"), "https://pan.baidu.com/s/1cIIZZv89eBKv255fKtD4cQ")

# Re-apply the hyperlink cell style, since Hyperlinks.Add() re-styles the
# cell with a fresh (non-shared) xf instead of reusing the existing one.
$ws.Range("G2").Style = "超链接"
$ws.Range("G3").Style = "超链接"
$ws.Range("G4").Style = "超链接"
$ws.Range("G5").Style = "超链接"
$ws.Range("G6").Style = "超链接"
$ws.Range("G7").Style = "超链接"
$ws.Range("G8").Style = "超链接"
$ws.Range("G9").Style = "超链接"
$ws.Range("G10").Style = "超链接"
$ws.Range("G11").Style = "超链接"

# ---------------------------------------------------------------------
# 4. Refresh the active selection to match the now-13-row table.
# ---------------------------------------------------------------------
$ws.Range("A10").Select()
